$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.837.28'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.749.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.01%  '

$ws.Range("E4").Value = '  +0.50%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("E6").Value = '  +0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3855'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3373'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.53'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07192'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.005'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.169'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.751.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.067'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001056'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06613'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '79.32'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.176'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.41%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.835.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.30%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.401'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.292'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.950.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.279'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '130.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.017'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.789'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08776'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.541'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6507'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.08%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.120'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02269'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06078'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2102'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.65%  '

$ws.Range("E42").Value = '  -2.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.018'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.48%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.817'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6013'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.74%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.989'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.31%  '

$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.151'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.18%  '

$ws.Range("B51").Value = 'Flow'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.096'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.96%  '
